$d = $word.ActiveDocument

# ------------------------------------------------------------------
# 1) Letter date: "September 19, 2025" -> "September 21, 2025"
# ------------------------------------------------------------------
$d.Content.Find.Execute(
    "September 19, 2025", $false, $false, $false, $false, $false,
    $true, 1, $false, "September 21, 2025", 2) | Out-Null

# ------------------------------------------------------------------
# 2) Split the mailing-address line into two paragraphs:
#      "121 9Th St, San Francisco CA 94103"
#    ->  "121 9Th St"
#        "San Francisco, CA 94103"
#    Only the FIRST occurrence (the mailing address near the top of
#    the letter) is touched; the later "PROPERTY ADDRESS:" field that
#    repeats the same string must stay untouched.
# ------------------------------------------------------------------
$addressPara = $null
foreach ($p in $d.Paragraphs) {
    $ptext = $p.Range.Text.TrimEnd([char]13)
    if ($ptext -eq "121 9Th St, San Francisco CA 94103") {
        $addressPara = $p
        break
    }
}

if ($addressPara -ne $null) {
    $r = $addressPara.Range
    $r.MoveEnd(1, -1) | Out-Null
    $r.Text = "121 9Th St"

    # Insert a brand-new paragraph right after it, inheriting the
    # same paragraph/run formatting, then fill in the city/state/zip.
    $addressPara.Range.InsertParagraphAfter() | Out-Null
    $cityPara = $addressPara.Next()
    $cr = $cityPara.Range
    $cr.MoveEnd(1, -1) | Out-Null
    $cr.Text = "San Francisco, CA 94103"
}

# ------------------------------------------------------------------
# 3) Remove the stray empty "No Spacing" paragraph that immediately
#    follows "... Board of Directors".
# ------------------------------------------------------------------
$boardPara = $null
foreach ($p in $d.Paragraphs) {
    $ptext = $p.Range.Text.TrimEnd([char]13)
    if ($ptext -eq "788 Minna Street Board of Directors") {
        $boardPara = $p
        break
    }
}

if ($boardPara -ne $null) {
    $trailingPara = $boardPara.Next()
    if ($trailingPara -ne $null -and $trailingPara.Range.Text.TrimEnd([char]13) -eq "") {
        $trailingPara.Range.Delete() | Out-Null
    }
}

"done"
